$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2675626666666667
$ws.Range("H2").Value = 0.802688
$ws.Range("I2").Value = 0.01337792263105189
$ws.Range("J2").Value = 0.01337792263105189
$ws.Range("M2").Value = 10.045207
$ws.Range("N2").Value = 30.135621
$ws.Range("O2").Value = 0.9365108453707793
$ws.Range("P2").Value = 0.9365108453707794
$ws.Range("Q2").Value = 2.687722372138666
$ws.Range("R2").Value = 24.189501349248
$ws.Range("S2").Value = 0.01252856963251129
$ws.Range("T2").Value = 0.01252856963251129
$ws.Range("G3").Value = 0.2675626666666667
$ws.Range("H3").Value = 0.802688
$ws.Range("I3").Value = 0.01337792263105189
$ws.Range("J3").Value = 0.01337792263105189
$ws.Range("O3").Value = 0.03971513502725754
$ws.Range("P3").Value = 0.03971513502725754
$ws.Range("Q3").Value = 0.1139797338737778
$ws.Range("R3").Value = 1.025817604864
$ws.Range("S3").Value = 0.0005313060036764304
$ws.Range("T3").Value = 0.0005313060036764304
$ws.Range("G4").Value = 0.2675626666666667
$ws.Range("H4").Value = 0.802688
$ws.Range("I4").Value = 0.01337792263105189
$ws.Range("J4").Value = 0.01337792263105189
$ws.Range("M4").Value = 0.255005
$ws.Range("N4").Value = 0.765015
$ws.Range("O4").Value = 0.02377401960196297
$ws.Range("P4").Value = 0.02377401960196297
$ws.Range("Q4").Value = 0.06822981781333333
$ws.Range("R4").Value = 0.61406836032
$ws.Range("S4").Value = 0.0003180469948641717
$ws.Range("T4").Value = 0.0003180469948641717
$ws.Range("I5").Value = 0.9475831922313891
$ws.Range("J5").Value = 0.947583192231389
$ws.Range("M5").Value = 10.045207
$ws.Range("N5").Value = 30.135621
$ws.Range("O5").Value = 0.9365108453707793
$ws.Range("P5").Value = 0.9365108453707794
$ws.Range("Q5").Value = 190.3763846945363
$ws.Range("R5").Value = 1713.387462250827
$ws.Range("S5").Value = 0.88742193641576
$ws.Range("T5").Value = 0.88742193641576
$ws.Range("I6").Value = 0.9475831922313891
$ws.Range("J6").Value = 0.947583192231389
$ws.Range("O6").Value = 0.03971513502725754
$ws.Range("P6").Value = 0.03971513502725754
$ws.Range("S6").Value = 0.03763339442902936
$ws.Range("T6").Value = 0.03763339442902935
$ws.Range("I7").Value = 0.9475831922313891
$ws.Range("J7").Value = 0.947583192231389
$ws.Range("M7").Value = 0.255005
$ws.Range("N7").Value = 0.765015
$ws.Range("O7").Value = 0.02377401960196297
$ws.Range("P7").Value = 0.02377401960196297
$ws.Range("Q7").Value = 4.832845154811666
$ws.Range("R7").Value = 43.495606393305
$ws.Range("S7").Value = 0.02252786138659969
$ws.Range("T7").Value = 0.02252786138659968
$ws.Range("A8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7737046666666667
$ws.Range("H8").Value = 2.321114
$ws.Range("I8").Value = 0.03868462405050454
$ws.Range("J8").Value = 0.03868462405050453
$ws.Range("M8").Value = 10.045207
$ws.Range("N8").Value = 30.135621
$ws.Range("O8").Value = 0.9365108453707793
$ws.Range("P8").Value = 0.9365108453707794
$ws.Range("Q8").Value = 7.772023533532667
$ws.Range("R8").Value = 69.948211801794
$ws.Range("S8").Value = 0.03622856997238879
$ws.Range("T8").Value = 0.03622856997238878
$ws.Range("A9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7737046666666667
$ws.Range("H9").Value = 2.321114
$ws.Range("I9").Value = 0.03868462405050454
$ws.Range("J9").Value = 0.03868462405050453
$ws.Range("O9").Value = 0.03971513502725754
$ws.Range("P9").Value = 0.03971513502725754
$ws.Range("Q9").Value = 0.3295925141657778
$ws.Range("R9").Value = 2.966332627492
$ws.Range("S9").Value = 0.001536365067644482
$ws.Range("T9").Value = 0.001536365067644482
$ws.Range("A10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7737046666666667
$ws.Range("H10").Value = 2.321114
$ws.Range("I10").Value = 0.03868462405050454
$ws.Range("J10").Value = 0.03868462405050453
$ws.Range("M10").Value = 0.255005
$ws.Range("N10").Value = 0.765015
$ws.Range("O10").Value = 0.02377401960196297
$ws.Range("P10").Value = 0.02377401960196297
$ws.Range("Q10").Value = 0.1972985585233333
$ws.Range("R10").Value = 1.77568702671
$ws.Range("S10").Value = 0.0009196890104712629
$ws.Range("T10").Value = 0.0009196890104712628
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.007085333333333333
$ws.Range("H11").Value = 0.021256
$ws.Range("I11").Value = 0.0003542610870545456
$ws.Range("J11").Value = 0.0003542610870545455
$ws.Range("M11").Value = 10.045207
$ws.Range("N11").Value = 30.135621
$ws.Range("O11").Value = 0.9365108453707793
$ws.Range("P11").Value = 0.9365108453707794
$ws.Range("Q11").Value = 0.07117363999733332
$ws.Range("R11").Value = 0.640562759976
$ws.Range("S11").Value = 0.0003317693501194237
$ws.Range("T11").Value = 0.0003317693501194237
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.007085333333333333
$ws.Range("H12").Value = 0.021256
$ws.Range("I12").Value = 0.0003542610870545456
$ws.Range("J12").Value = 0.0003542610870545455
$ws.Range("O12").Value = 0.03971513502725754
$ws.Range("P12").Value = 0.03971513502725754
$ws.Range("Q12").Value = 0.003018300040888889
$ws.Range("R12").Value = 0.027164700368
$ws.Range("S12").Value = 0.00001406952690727432
$ws.Range("T12").Value = 0.00001406952690727431
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.007085333333333333
$ws.Range("H13").Value = 0.021256
$ws.Range("I13").Value = 0.0003542610870545456
$ws.Range("J13").Value = 0.0003542610870545455
$ws.Range("M13").Value = 0.255005
$ws.Range("N13").Value = 0.765015
$ws.Range("O13").Value = 0.02377401960196297
$ws.Range("P13").Value = 0.02377401960196297
$ws.Range("Q13").Value = 0.001806795426666667
$ws.Range("R13").Value = 0.01626115884
$ws.Range("S13").Value = 0.000008422210027847475
$ws.Range("T13").Value = 0.000008422210027847473
